$d = $word.ActiveDocument

# 1) Collapse the "Easy Save" sentence (previously split across three runs
#    around a spell-check proofErr bracket for "Easy") into a single run.
$d.Content.Find.Execute(
    "Lorsque vous lancez l’application Easy Save pour la première fois, une interface vous demandant quel langage utiliser apparait. Le langage choisi est configuré dans un fichier et peut être modifié.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Lorsque vous lancez l’application Easy Save pour la première fois, une interface vous demandant quel langage utiliser apparait. Le langage choisi est configuré dans un fichier et peut être modifié.",
    2) | Out-Null

# 2) Add a new bulleted item "Quitter l'application" right after
#    "Obtenir les informations d'un projet", matching the existing list
#    formatting (numId 1 / ilvl 0).
$rng = $d.Content
$rng.Find.Execute("Obtenir les informations d’un projet") | Out-Null
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$newPara = $rng.Paragraphs(1).Next()
$newPara.Range.Text = "Quitter l’application"

# 3) Collapse the "répertoire Debug du projet" sentence (previously split
#    across three runs around a spell-check proofErr bracket for "Debug")
#    into a single run.
$d.Content.Find.Execute(
    "Lorsque le projet de sauvegarde se lance, les informations sur ce qu’il se passe sont stockées dans un dossier précis dans le répertoire Debug du projet.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Lorsque le projet de sauvegarde se lance, les informations sur ce qu’il se passe sont stockées dans un dossier précis dans le répertoire Debug du projet.",
    2) | Out-Null
